# Generate Report for Handoff
# Updates the "b.md" row across all three sheets to reflect that a new
# handoff (b.*.xlf) went out, but the handback version is stale.

$wb = $excel.ActiveWorkbook

# ---- Overview sheet ---------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = "Ready for handoff"
$wsOverview.Range("F3").Value = "Ready for handoff"
$wsOverview.Range("G3").Value = "2016-09-01 00:40:42"

# ---- zh-cn sheet --------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = "Ready for handoff"
$wsZhCn.Range("F3").Value = "'False"
$wsZhCn.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$wsZhCn.Range("H3").Value = "2016-09-01 00:40:37"
$wsZhCn.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/3f18b0e0e40a697767e0879ce7f04eac0c03d510/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/a0fcbcbfef932390ec8bd9c42090ad81e5daf807/e2e/b.md."
$wsZhCn.Columns.Item(16).ColumnWidth = 39.1

# ---- de-de sheet --------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = "Ready for handoff"
$wsDeDe.Range("F3").Value = "False"
$wsDeDe.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$wsDeDe.Range("H3").Value = "2016-09-01 00:40:42"
$wsDeDe.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/3f18b0e0e40a697767e0879ce7f04eac0c03d510/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/a0fcbcbfef932390ec8bd9c42090ad81e5daf807/e2e/b.md."
$wsDeDe.Columns.Item(16).ColumnWidth = 39.1
